$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88, pushing the existing row 88 down to row 89.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the data that used to be in row 87
# (a duplicate weekly entry for the same period, prior to row 87 being updated).
$ws.Range("A88").Value = 11
$ws.Range("B88").Value = "Vega Monumental Concepción"
$ws.Range("C88").Value = "Bíobío"
$ws.Range("D88").Value = 44664
$ws.Range("E88").Value = 8
$ws.Range("F88").Value = 100112001
$ws.Range("G88").Value = "Berenjena"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 100
$ws.Range("K88").Value = 6500
$ws.Range("L88").Value = 7000
$ws.Range("M88").Value = 6750
$ws.Range("N88").Value = "$/caja 60 unidades"
$ws.Range("O88").Value = "Región de Arica y Parinacota"
$ws.Range("P88").Value = 112
$ws.Range("Q88").Value = 60
$ws.Range("R88").Value = "Hortaliza"

# Update row 87 with the new weekly values.
$ws.Range("D87").Value = 44706
$ws.Range("K87").Value = 5500
$ws.Range("L87").Value = 6000
$ws.Range("M87").Value = 5750
$ws.Range("P87").Value = 96
